$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: find the row number whose column A matches a given label exactly
function Find-RowByLabel($label) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value2 -eq $label) {
            return $r
        }
    }
    return -1
}

# 1) Delete the entire row for "RM 232"
$rowRM232 = Find-RowByLabel "RM 232"
if ($rowRM232 -gt 0) {
    $ws.Rows.Item($rowRM232).Delete()
}

# 2) Delete the entire row for "SC 92" (rows have shifted up after step 1)
$rowSC92 = Find-RowByLabel "SC 92"
if ($rowSC92 -gt 0) {
    $ws.Rows.Item($rowSC92).Delete()
}

# 3) Update the "SC 5" row: column D (value C) becomes -13.8
$rowSC5 = Find-RowByLabel "SC 5"
if ($rowSC5 -gt 0) {
    $ws.Cells.Item($rowSC5, 4).Value = -13.8
}

# 4) Clear the "SC 101" row: column D (value C) becomes empty
$rowSC101 = Find-RowByLabel "SC 101"
if ($rowSC101 -gt 0) {
    $ws.Cells.Item($rowSC101, 4).ClearContents()
}

# 5) Update the "SC 232" row: column E (value D) becomes -10.7
$rowSC232 = Find-RowByLabel "SC 232"
if ($rowSC232 -gt 0) {
    $ws.Cells.Item($rowSC232, 5).Value = -10.7
}
